# Update countries & provincias Spain
#
# 1. Swap the displayed country names for the two tied rows (both had
#    "Casos totales" = 26): row 202 becomes "Santa Lucia", row 203 becomes
#    "Timor Oriental".
# 2. Bump the "Datos actualizados" timestamp from 07:03 to 08:20.
# 3. Refresh the case-count figures for five country rows (33, 62, 72, 73,
#    150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Santa Lucia / Timor Oriental swap (rows 202-203, column A) ---
$ws.Range("A202").Value2 = "Santa Lucia"
$ws.Range("A203").Value2 = "Timor Oriental"

# --- 2. Timestamp update (row 1) ---
$ws.Range("A1").Value2 = "Datos actualizados a 24 de Agosto de 2020 a las 08:20"

# --- 3. Numeric refresh ---
# Row 33 - Israel
$ws.Range("B33").Value2 = 103151
$ws.Range("C33").Value2 = 488
$ws.Range("D33").Value2 = 80511
$ws.Range("E33").Value2 = 21806

# Row 62 - Uzbekistan
$ws.Range("B62").Value2 = 39065
$ws.Range("C62").Value2 = 119
$ws.Range("E62").Value2 = 3803
$ws.Range("G62").Value2 = 2
$ws.Range("H62").Value2 = 275

# Row 72 - Australia
$ws.Range("B72").Value2 = 24916
$ws.Range("C72").Value2 = 104
$ws.Range("D72").Value2 = 19234

# Row 73 - El Salvador
$ws.Range("D73").Value2 = 12276
$ws.Range("E73").Value2 = 11677
$ws.Range("G73").Value2 = 8
$ws.Range("H73").Value2 = 669

# Row 150 - Georgia
$ws.Range("B150").Value2 = 1421
$ws.Range("C150").Value2 = 10
$ws.Range("D150").Value2 = 1137
$ws.Range("E150").Value2 = 267
